$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.548388957977295
$ws.Range("B1").Value = 2.472123384475708
$ws.Range("C1").Value = 1.894278526306152
$ws.Range("D1").Value = 1.703905701637268
$ws.Range("E1").Value = 1.530171751976013
